# Adds a new "Changes in inventories" indicator row to the Cataloged_Indicators
# sheet (inserted right above the "9. Consumer" section, which pushes that
# section's rows down by one), and updates the sheet's view state to reflect
# where the author was scrolled to / what was selected when the edit was made.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 68 -- this shifts the old row 68 (blank
# separator) and the whole "9. Consumer" block (old rows 69-73) down to
# rows 69-74, exactly like the diff shows.
$ws.Rows.Item(68).Insert()

# Fill the newly inserted row 68 with the new indicator's data.
$ws.Range("B68").Value = "Changes in inventories"
$ws.Range("C68").Value = "OECD"
$ws.Range("D68").Value = "Q"
$ws.Range("E68").Value = "Current prices, local currency"
$ws.Range("F68").Value = 80

# Reflect the author's final view state: still frozen on row 1 only, but
# scrolled further down, with cell B69 selected.
$ws.Activate()
$ws.Range("B69").Select()
